{"js": "const replacements = [\n  [\"2024-02-19 Monday\", \"2024-02-20 Tuesday\"],\n  [\"75\u00f72=37, 1\", \"27\u00f73=9, 0\"],\n  [\"71\u00f72=35, 1\", \"72\u00f79=8, 0\"],\n  [\"71\u00f78=8, 7\", \"65\u00f78=8, 1\"],\n  [\"98\u00f78=12, 2\", \"61\u00f75=12, 1\"],\n  [\"68\u00f74=17, 0\", \"14\u00f73=4, 2\"],\n  [\"56\u00f77=8, 0\", \"43\u00f79=4, 7\"],\n  [\"62\u00f73=20, 2\", \"99\u00f78=12, 3\"],\n  [\"41\u00f79=4, 5\", \"45\u00f73=15, 0\"],\n  [\"97\u00f75=19, 2\", \"75\u00f75=15, 0\"],\n  [\"47\u00f73=15, 2\", \"90\u00f75=18, 0\"],\n  [\"23\u00f72=11, 1\", \"54\u00f73=18, 0\"],\n  [\"71\u00f77=10, 1\", \"25\u00f74=6, 1\"],\n  [\"46\u00f78=5, 6\", \"56\u00f76=9, 2\"],\n  [\"57\u00f73=19, 0\", \"74\u00f78=9, 2\"],\n  [\"18\u00f74=4, 2\", \"53\u00f77=7, 4\"],\n  [\"76\u00f77=10, 6\", \"27\u00f79=3, 0\"],\n  [\"84\u00f76=14, 0\", \"92\u00f78=11, 4\"],\n  [\"98\u00f72=49, 0\", \"46\u00f75=9, 1\"],\n  [\"80\u00f77=11, 3\", \"93\u00f74=23, 1\"],\n  [\"73\u00f74=18, 1\", \"89\u00f75=17, 4\"],\n  [\"54\u00f72=27, 0\", \"80\u00f74=20, 0\"],\n  [\"21\u00f73=7, 0\", \"67\u00f72=33, 1\"],\n  [\"24\u00f79=2, 6\", \"36\u00f77=5, 1\"],\n  [\"44\u00f78=5, 4\", \"18\u00f78=2, 2\"],\n  [\"92\u00f79=10, 2\", \"14\u00f73=4, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-02-19 Monday\", \"2024-02-20 Tuesday\"),\n  @(\"75\u00f72=37, 1\", \"27\u00f73=9, 0\"),\n  @(\"71\u00f72=35, 1\", \"72\u00f79=8, 0\"),\n  @(\"71\u00f78=8, 7\", \"65\u00f78=8, 1\"),\n  @(\"98\u00f78=12, 2\", \"61\u00f75=12, 1\"),\n  @(\"68\u00f74=17, 0\", \"14\u00f73=4, 2\"),\n  @(\"56\u00f77=8, 0\", \"43\u00f79=4, 7\"),\n  @(\"62\u00f73=20, 2\", \"99\u00f78=12, 3\"),\n  @(\"41\u00f79=4, 5\", \"45\u00f73=15, 0\"),\n  @(\"97\u00f75=19, 2\", \"75\u00f75=15, 0\"),\n  @(\"47\u00f73=15, 2\", \"90\u00f75=18, 0\"),\n  @(\"23\u00f72=11, 1\", \"54\u00f73=18, 0\"),\n  @(\"71\u00f77=10, 1\", \"25\u00f74=6, 1\"),\n  @(\"46\u00f78=5, 6\", \"56\u00f76=9, 2\"),\n  @(\"57\u00f73=19, 0\", \"74\u00f78=9, 2\"),\n  @(\"18\u00f74=4, 2\", \"53\u00f77=7, 4\"),\n  @(\"76\u00f77=10, 6\", \"27\u00f79=3, 0\"),\n  @(\"84\u00f76=14, 0\", \"92\u00f78=11, 4\"),\n  @(\"98\u00f72=49, 0\", \"46\u00f75=9, 1\"),\n  @(\"80\u00f77=11, 3\", \"93\u00f74=23, 1\"),\n  @(\"73\u00f74=18, 1\", \"89\u00f75=17, 4\"),\n  @(\"54\u00f72=27, 0\", \"80\u00f74=20, 0\"),\n  @(\"21\u00f73=7, 0\", \"67\u00f72=33, 1\"),\n  @(\"24\u00f79=2, 6\", \"36\u00f77=5, 1\"),\n  @(\"44\u00f78=5, 4\", \"18\u00f78=2, 2\"),\n  @(\"92\u00f79=10, 2\", \"14\u00f73=4, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
